# Inventory export re-run: rename the inventory sheet to the new count date,
# update the A0001 line, and append the new A0002 recount line.
#
# NOTE: the PNR/Count columns in this export are stored as *text* (not
# numbers) even though their contents look numeric ("2000", "5", "2766",
# "8"). Plain `Range.Value = "2000"` gets auto-coerced to a number by the
# COM layer (same as typing 2000 into a General-formatted cell in real
# Excel), so we briefly mark the cell as Text ("@"), write the value, and
# then paste back the Normal/General formatting from A1 (which keeps the
# worksheet's original default style index) while leaving the stored cell
# type as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range("A1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats - restore default style/look
}

# Rename sheet: "Inventúra 2023-12-23" -> "Inventúra 2024-01-15"
$ws.Name = "Inventúra 2024-01-15"

# Row 2 (A0001): PNR/Material/Count updated
Set-TextValue "B2" "2000"
$ws.Range("C2").Value = "Test material"
Set-TextValue "D2" "5"

# Row 3: was A0002/002/Skrina/2 -> now a second A0001 line
$ws.Range("A3").Value = "A0001"
Set-TextValue "B3" "2766"
$ws.Range("C3").Value = "ds"
Set-TextValue "D3" "8"

# Row 4 (new): A0002 recount line, same PNR/Material/Count as row 2
$ws.Range("A4").Value = "A0002"
Set-TextValue "B4" "2000"
$ws.Range("C4").Value = "Test material"
Set-TextValue "D4" "5"
